$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date on the existing accident record (row 2, column G)
$ws.Range("G2").Value = "17/11/2019"

# Add new accident record - row 3
$ws.Range("A3").Value = "8022"
$ws.Range("B3").Value = "Cao Thành Duy"
$ws.Range("C3").Value = "Cơ điện vận tải"
$ws.Range("D3").Value = "dsgfsg"
$ws.Range("E3").Value = "Nặng"
$ws.Range("F3").Value = "CA 1"
$ws.Range("G3").Value = "25/10/2019"

# Add new accident record - row 4
$ws.Range("A4").Value = "8022"
$ws.Range("B4").Value = "Cao Thành Duy"
$ws.Range("C4").Value = "Cơ điện vận tải"
$ws.Range("D4").Value = "fgdh"
$ws.Range("E4").Value = "Nặng"
$ws.Range("F4").Value = "CA 1"
$ws.Range("G4").Value = "15/10/2019"

# Move the active selection to G4, matching the saved workbook state
$ws.Range("G4").Select()
